$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were inserted before the old row 232, pushing the
# existing rows 232-243 down to 234-245 (dimension grows from R243 to R245).
$ws.Rows.Item(232).Insert()
$ws.Rows.Item(232).Insert()

# New row 232 (Pepino dulce, Primera, week of 2022-07-04)
$ws.Range("A232").Value = 10
$ws.Range("B232").Value = "Vega Modelo de Temuco"
$ws.Range("C232").Value = "La Araucanía"
$ws.Range("D232").Value = 44746
$ws.Range("E232").Value = 9
$ws.Range("F232").Value = 100112043
$ws.Range("G232").Value = "Pepino dulce"
$ws.Range("H232").Value = "Cultivar IV Región"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 110
$ws.Range("K232").Value = 15000
$ws.Range("L232").Value = 15000
$ws.Range("M232").Value = 15000
$ws.Range("N232").Value = "$/bandeja 18 kilos"
$ws.Range("O232").Value = "Provincia de Limarí"
$ws.Range("P232").Value = 833
$ws.Range("Q232").Value = 18
$ws.Range("R232").Value = "Hortaliza"

# New row 233 (Pepino dulce, Segunda, week of 2022-07-04)
$ws.Range("A233").Value = 10
$ws.Range("B233").Value = "Vega Modelo de Temuco"
$ws.Range("C233").Value = "La Araucanía"
$ws.Range("D233").Value = 44746
$ws.Range("E233").Value = 9
$ws.Range("F233").Value = 100112043
$ws.Range("G233").Value = "Pepino dulce"
$ws.Range("H233").Value = "Cultivar IV Región"
$ws.Range("I233").Value = "Segunda"
$ws.Range("J233").Value = 65
$ws.Range("K233").Value = 12000
$ws.Range("L233").Value = 12000
$ws.Range("M233").Value = 12000
$ws.Range("N233").Value = "$/bandeja 18 kilos"
$ws.Range("O233").Value = "Provincia de Limarí"
$ws.Range("P233").Value = 667
$ws.Range("Q233").Value = 18
$ws.Range("R233").Value = "Hortaliza"
